# "Le he mostrado a Marta como funciona el GIT"
# Add a 3rd slide (same "Title and Content" layout as slide 2) with the
# title "Hola marta" and an empty content placeholder.

$p = $ppt.ActivePresentation

# Duplicate slide 2 (Título y objetos / Title-and-Content layout) so the
# new slide 3 inherits the exact same placeholder/paragraph structure
# (incl. the group xfrm and the content placeholder's empty endParaRPr)
# instead of the bare structure produced by Slides.Add.
$s2 = $p.Slides.Item(2)
$newSlide = $s2.Duplicate()

$titleRange = $newSlide.Shapes.Item(1).TextFrame.TextRange

# The duplicated title text is "Una " / "nueva" / " " / "diapositiva"
# (4 runs). Rewrite each run in place (instead of replacing the whole
# string) so the run boundaries - and which run keeps a dirty="0" marker
# versus not - end up matching "Hola " / "marta".
$run4 = $titleRange.Characters(11, 11)
$run4.Text = "marta"

$run3 = $titleRange.Characters(10, 1)
$run3.Text = ""

$run2 = $titleRange.Characters(5, 5)
$run2.Text = ""

$run1 = $titleRange.Characters(1, 4)
$run1.Text = "Hola "

# Content placeholder (shape 2) stays empty, exactly as duplicated.
